# Generate Report for Handback
#
# For the "3c9d52dc-9345-4bfa-bf4d-ae398e48c57a" row (row 7) on both the
# zh-cn and de-de sheets, a handback xliff was generated but its version
# is stale compared to the latest handoff. Record that in the report:
#   - Latest Target File   (I7) -> link to the source .md file (same
#                                   hyperlink target as column A)
#   - Latest Handback File (J7) -> the generated .xlf file name
#   - Latest Handback DateTime (K7) -> when it was generated
#   - Error Detail (P7)    -> explanation that the handback is stale
# Also widen the Error Detail column so the message is readable.

$wb = $excel.ActiveWorkbook

$mdFile = "3c9d52dc-9345-4bfa-bf4d-ae398e48c57a.md"
$mdUrl  = "https://github.com/OpenLocalizationTestOrg/oltest/blob/4cb3291eb9d2fcb11a97872cca35cf954645e308/e2e/3c9d52dc-9345-4bfa-bf4d-ae398e48c57a.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/7f4e1aea1f19941fc3806d739c754fa9a8b7d32a/e2e/3c9d52dc-9345-4bfa-bf4d-ae398e48c57a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/4cb3291eb9d2fcb11a97872cca35cf954645e308/e2e/3c9d52dc-9345-4bfa-bf4d-ae398e48c57a.md."

function Update-LocalizationSheet {
    param(
        $ws,
        [string]$targetXlf,
        [string]$handbackDateTime
    )

    # Latest Target File: hyperlink to the markdown source, same as A7.
    $ws.Range("I7").Value = $mdFile
    $ws.Hyperlinks.Add($ws.Range("I7"), $mdUrl, "", "", $mdFile) | Out-Null
    $ws.Range("I7").Font.Underline = 2
    $ws.Range("I7").Font.Color = 15570276

    # Latest Handback File / DateTime.
    $ws.Range("J7").Value = $targetXlf
    $ws.Range("K7").Value = $handbackDateTime

    # Error Detail.
    $ws.Range("P7").Value = $errorDetail

    # Widen the Error Detail column so the message is legible.
    $ws.Columns.Item(16).ColumnWidth = 39.14
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LocalizationSheet -ws $wsZhCn `
    -targetXlf "3c9d52dc-9345-4bfa-bf4d-ae398e48c57a.e4e5273c8c8643fd03cf6153fee2cf8a0ba98300.zh-cn.xlf" `
    -handbackDateTime "2016-08-13 14:52:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LocalizationSheet -ws $wsDeDe `
    -targetXlf "3c9d52dc-9345-4bfa-bf4d-ae398e48c57a.e4e5273c8c8643fd03cf6153fee2cf8a0ba98300.de-de.xlf" `
    -handbackDateTime "2016-08-13 14:53:03"

Write-Output "Updated handback report for 3c9d52dc-9345-4bfa-bf4d-ae398e48c57a on zh-cn and de-de sheets"
